$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$blockA = New-Object 'object[,]' 24,5
$blockA[0,0] = 1.02
$blockA[0,1] = 1.041992806374259
$blockA[0,2] = 1.042311102517584
$blockA[0,3] = 1.04950866930966
$blockA[0,4] = 1.058736063525053
$blockA[1,0] = 1.02
$blockA[1,1] = 1.043148121805795
$blockA[1,2] = 1.043321539733336
$blockA[1,3] = 1.050540462959215
$blockA[1,4] = 1.059867459277227
$blockA[2,0] = 1.019999999999999
$blockA[2,1] = 1.043895304716095
$blockA[2,2] = 1.043975288785779
$blockA[2,3] = 1.051208106379966
$blockA[2,4] = 1.060599546329461
$blockA[3,0] = 1.02
$blockA[3,1] = 1.044209329887426
$blockA[3,2] = 1.044250108266513
$blockA[3,3] = 1.051488785098123
$blockA[3,4] = 1.060907315967653
$blockA[4,0] = 1.02
$blockA[4,1] = 1.044262050823873
$blockA[4,2] = 1.044296250705476
$blockA[4,3] = 1.051535912387333
$blockA[4,4] = 1.060958991876702
$blockA[5,0] = 1.02
$blockA[5,1] = 1.043899501089927
$blockA[5,2] = 1.043978961001867
$blockA[5,3] = 1.051211856816553
$blockA[5,4] = 1.060603658760702
$blockA[6,0] = 1.02
$blockA[6,1] = 1.042383330790915
$blockA[6,2] = 1.042652599573664
$blockA[6,3] = 1.049857368282253
$blockA[6,4] = 1.05911842474552
$blockA[7,0] = 1.02
$blockA[7,1] = 1.039708650855793
$blockA[7,2] = 1.04031481252085
$blockA[7,3] = 1.047470594093767
$blockA[7,4] = 1.056501210507586
$blockA[8,0] = 1.02
$blockA[8,1] = 1.037923434120079
$blockA[8,2] = 1.038755871246676
$blockA[8,3] = 1.045879381574178
$blockA[8,4] = 1.054756327795801
$blockA[9,0] = 1.02
$blockA[9,1] = 1.037149898622302
$blockA[9,2] = 1.038080724989075
$blockA[9,3] = 1.045190353106145
$blockA[9,4] = 1.054000745829484
$blockA[10,0] = 1.02
$blockA[10,1] = 1.036862492432356
$blockA[10,2] = 1.037829927557086
$blockA[10,3] = 1.044934412929149
$blockA[10,4] = 1.053720082517481
$blockA[11,0] = 1.02
$blockA[11,1] = 1.0369241457151
$blockA[11,2] = 1.037883725288425
$blockA[11,3] = 1.044989313175833
$blockA[11,4] = 1.053780286063524
$blockA[12,0] = 1.02
$blockA[12,1] = 1.037126143193009
$blockA[12,2] = 1.038059994368378
$blockA[12,3] = 1.045169197104612
$blockA[12,4] = 1.053977546252021
$blockA[13,0] = 1.02
$blockA[13,1] = 1.037250589762831
$blockA[13,2] = 1.038168597148205
$blockA[13,3] = 1.045280028930539
$blockA[13,4] = 1.054099083867038
$blockA[14,0] = 1.02
$blockA[14,1] = 1.037974759910918
$blockA[14,2] = 1.038800675996491
$blockA[14,3] = 1.045925109566202
$blockA[14,4] = 1.054806472426758
$blockA[15,0] = 1.02
$blockA[15,1] = 1.038428870916637
$blockA[15,2] = 1.039197131124644
$blockA[15,3] = 1.04632974484669
$blockA[15,4] = 1.05525018811862
$blockA[16,0] = 1.02
$blockA[16,1] = 1.038693695662799
$blockA[16,2] = 1.03942836571415
$blockA[16,3] = 1.046565759605325
$blockA[16,4] = 1.055508996551471
$blockA[17,0] = 1.02
$blockA[17,1] = 1.03878398554326
$blockA[17,2] = 1.039507208897384
$blockA[17,3] = 1.046646234199933
$blockA[17,4] = 1.055597242991021
$blockA[18,0] = 1.02
$blockA[18,1] = 1.038380154341733
$blockA[18,2] = 1.039154596383066
$blockA[18,3] = 1.046286331539503
$blockA[18,4] = 1.055202581985815
$blockA[19,0] = 1.02
$blockA[19,1] = 1.037066662211723
$blockA[19,2] = 1.038008088023007
$blockA[19,3] = 1.045116225897561
$blockA[19,4] = 1.053919458232028
$blockA[20,0] = 1.02
$blockA[20,1] = 1.036240349975631
$blockA[20,2] = 1.037287127839744
$blockA[20,3] = 1.044380508574765
$blockA[20,4] = 1.053112669091582
$blockA[21,0] = 1.02
$blockA[21,1] = 1.036678438656171
$blockA[21,2] = 1.037669332705825
$blockA[21,3] = 1.044770528957219
$blockA[21,4] = 1.053540367275746
$blockA[22,0] = 1.02
$blockA[22,1] = 1.038402167420527
$blockA[22,2] = 1.039173816035497
$blockA[22,3] = 1.046305948150645
$blockA[22,4] = 1.055224093158242
$blockA[23,0] = 1.02
$blockA[23,1] = 1.040400483134435
$blockA[23,2] = 1.040919257391312
$blockA[23,3] = 1.048087635078025
$blockA[23,4] = 1.057177832698101
$ws.Range("B2:F25").Value = $blockA

$blockB = New-Object 'object[,]' 24,6
$blockB[0,0] = 1.039504644312983
$blockB[0,1] = 1.047071025778354
$blockB[0,2] = 1.045088148838613
$blockB[0,3] = 1.052265531255145
$blockB[0,4] = 1.061467483666029
$blockB[0,5] = 1.019547718600481
$blockB[1,0] = 1.039750573599286
$blockB[1,1] = 1.047871573272498
$blockB[1,2] = 1.045908868592163
$blockB[1,3] = 1.053109018329469
$blockB[1,4] = 1.062412167181153
$blockB[1,5] = 1.019816214211049
$blockB[2,0] = 1.039907633106465
$blockB[2,1] = 1.048388663739822
$blockB[2,2] = 1.046439250706199
$blockB[2,3] = 1.053654210092169
$blockB[2,4] = 1.063022852724481
$blockB[2,5] = 1.019989555302213
$blockB[3,0] = 1.039973164601963
$blockB[3,1] = 1.048605829480958
$blockB[3,2] = 1.046662061395726
$blockB[3,3] = 1.053883265356826
$blockB[3,4] = 1.06327944472177
$blockB[3,5] = 1.020062333795239
$blockB[4,0] = 1.039984138523851
$blockB[4,1] = 1.048642279731539
$blockB[4,2] = 1.046699462807549
$blockB[4,3] = 1.053921716351071
$blockB[4,4] = 1.063322519440729
$blockB[4,5] = 1.020074548109396
$blockB[5,0] = 1.039908510691732
$blockB[5,1] = 1.048391566379345
$blockB[5,2] = 1.046442228549326
$blockB[5,3] = 1.053657271302527
$blockB[5,4] = 1.063026281869575
$blockB[5,5] = 1.019990528141625
$blockB[6,0] = 1.03958818653827
$blockB[6,1] = 1.047341764793549
$blockB[6,2] = 1.045365655684404
$blockB[6,3] = 1.052550716344243
$blockB[6,4] = 1.061786865925042
$blockB[6,5] = 1.019638539496811
$blockB[7,0] = 1.03900785543467
$blockB[7,1] = 1.045484831204758
$blockB[7,2] = 1.043463378773376
$blockB[7,3] = 1.05059620181176
$blockB[7,4] = 1.059598329728267
$blockB[7,5] = 1.019015271457911
$blockB[8,0] = 1.038610288404494
$blockB[8,1] = 1.04424209984683
$blockB[8,2] = 1.042191652228365
$blockB[8,3] = 1.049290049698313
$blockB[8,4] = 1.058136224047604
$blockB[8,5] = 1.018597721231473
$blockB[9,0] = 1.038435602952328
$blockB[9,1] = 1.043702841706964
$blockB[9,2] = 1.041640132051248
$blockB[9,3] = 1.048723717555012
$blockB[9,4] = 1.057502375750612
$blockB[9,5] = 1.018416431730334
$blockB[10,0] = 1.038370335725503
$blockB[10,1] = 1.043502363968446
$blockB[10,2] = 1.041435143730343
$blockB[10,3] = 1.048513241606165
$blockB[10,4] = 1.05726682311943
$blockB[10,5] = 1.018349019264007
$blockB[11,0] = 1.038384353024201
$blockB[11,1] = 1.043545374978467
$blockB[11,2] = 1.041479120279827
$blockB[11,3] = 1.048558394628186
$blockB[11,4] = 1.05731735509397
$blockB[11,5] = 1.01836348279548
$blockB[12,0] = 1.038430215722137
$blockB[12,1] = 1.043686273694099
$blockB[12,2] = 1.041623190292619
$blockB[12,3] = 1.048706321897794
$blockB[12,4] = 1.057482907210369
$blockB[12,5] = 1.018410860895036
$blockB[13,0] = 1.038458422715017
$blockB[13,1] = 1.043773063058248
$blockB[13,2] = 1.041711939440815
$blockB[13,3] = 1.048797449510106
$blockB[13,4] = 1.05758489431122
$blockB[13,5] = 1.018440042362859
$blockB[14,0] = 1.038621828246933
$blockB[14,1] = 1.044277864386185
$blockB[14,2] = 1.042228236737977
$blockB[14,3] = 1.049327619220458
$blockB[14,4] = 1.058178274655456
$blockB[14,5] = 1.018609742527796
$blockB[15,0] = 1.038723649025848
$blockB[15,1] = 1.044594205141945
$blockB[15,2] = 1.042551867112498
$blockB[15,3] = 1.049659976697514
$blockB[15,4] = 1.058550285845212
$blockB[15,5] = 1.018716060277721
$blockB[16,0] = 1.038782794662157
$blockB[16,1] = 1.044778610719032
$blockB[16,2] = 1.042740552874267
$blockB[16,3] = 1.049853761951955
$blockB[16,4] = 1.058767201677963
$blockB[16,5] = 1.018778026601828
$blockB[17,0] = 1.038802920296688
$blockB[17,1] = 1.044841469516942
$blockB[17,2] = 1.042804875898571
$blockB[17,3] = 1.049919825318828
$blockB[17,4] = 1.058841152187257
$blockB[17,5] = 1.018799147548503
$blockB[18,0] = 1.038712749931987
$blockB[18,1] = 1.044560276222266
$blockB[18,2] = 1.042517153160884
$blockB[18,3] = 1.049624325466966
$blockB[18,4] = 1.058510380019123
$blockB[18,5] = 1.018704658254517
$blockB[19,0] = 1.038416720827477
$blockB[19,1] = 1.043644787324046
$blockB[19,2] = 1.041580768837222
$blockB[19,3] = 1.048662764184582
$blockB[19,4] = 1.057434159369503
$blockB[19,5] = 1.018396911255703
$blockB[20,0] = 1.038228389543034
$blockB[20,1] = 1.043068180270264
$blockB[20,2] = 1.040991279084689
$blockB[20,3] = 1.048057526808205
$blockB[20,4] = 1.05675684159922
$blockB[20,5] = 1.018202993289584
$blockB[21,0] = 1.038328436702887
$blockB[21,1] = 1.043373945938074
$blockB[21,2] = 1.041303849917729
$blockB[21,3] = 1.048378437971846
$blockB[21,4] = 1.057115962980214
$blockB[21,5] = 1.018305833256347
$blockB[22,0] = 1.038717675519301
$blockB[22,1] = 1.044575607582131
$blockB[22,2] = 1.042532839155543
$blockB[22,3] = 1.049640434949445
$blockB[22,4] = 1.058528411966328
$blockB[22,5] = 1.018709810483197
$blockB[23,0] = 1.039159767007103
$blockB[23,1] = 1.045965732092436
$blockB[23,2] = 1.043955784739794
$blockB[23,3] = 1.051102042017098
$blockB[23,4] = 1.06016465912204
$blockB[23,5] = 1.019176759808525
$ws.Range("I2:N25").Value = $blockB
